# Further adjusting of matlab scripts for reading data
#
# The "Road" worksheet has two blank rows removed:
#   - the blank row between the header row and the "Mean" row
#   - the blank row between the "STD DEV" row and the "Max" row
# Excel automatically re-numbers the remaining rows and rewrites the
# formulas that reference them when EntireRow is deleted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Road")

# Delete the empty row right below the header (old row 2)
$ws.Rows(2).Delete()

# Delete the empty row that is now between "STD DEV" and "Max" (old row 5,
# now row 4 after the previous deletion shifted everything up by one)
$ws.Rows(4).Delete()

# Select the two remaining empty rows (now rows 6 and 7), matching the
# selection state left behind by the edit.
$ws.Range("A6:A7").EntireRow.Select()
